$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Insert new record near the top of the Berenjena block (new row 205) ---
$ws.Rows(205).Insert()
$ws.Range("A205").Value = 8
$ws.Range("B205").Value = "Terminal La Palmera de La Serena"
$ws.Range("C205").Value = "Coquimbo"
$ws.Range("D205").Value = 45120
$ws.Range("E205").Value = 4
$ws.Range("F205").Value = 100112001
$ws.Range("G205").Value = "Berenjena"
$ws.Range("H205").Value = "Sin especificar"
$ws.Range("I205").Value = "Primera"
$ws.Range("J205").Value = 340
$ws.Range("K205").Value = 8000
$ws.Range("L205").Value = 8500
$ws.Range("M205").Value = 8250
$ws.Range("N205").Value = "$/caja 50 unidades"
$ws.Range("O205").Value = "Región de Arica y Parinacota"
$ws.Range("P205").Value = 165
$ws.Range("Q205").Value = 50
$ws.Range("R205").Value = "Hortaliza"

# --- Insert a second new record near the end of the block (new row 265) ---
$ws.Rows(265).Insert()
$ws.Range("A265").Value = 8
$ws.Range("B265").Value = "Terminal La Palmera de La Serena"
$ws.Range("C265").Value = "Coquimbo"
$ws.Range("D265").Value = 45121
$ws.Range("E265").Value = 4
$ws.Range("F265").Value = 100112001
$ws.Range("G265").Value = "Berenjena"
$ws.Range("H265").Value = "Sin especificar"
$ws.Range("I265").Value = "Primera"
$ws.Range("J265").Value = 500
$ws.Range("K265").Value = 8000
$ws.Range("L265").Value = 9000
$ws.Range("M265").Value = 8500
$ws.Range("N265").Value = "$/caja 50 unidades"
$ws.Range("O265").Value = "Región de Arica y Parinacota"
$ws.Range("P265").Value = 170
$ws.Range("Q265").Value = 50
$ws.Range("R265").Value = "Hortaliza"

Write-Output "done"
